# ---------------------------------------------------------------------------
# Update "2._OEM_Participant.xlsx" (Eq Guinea OEM Participant XLSForm)
#   - survey sheet: split the single "p_birth_date" date question into a
#     "birth_data" group with three integer sub-questions
#     (p_bd_yrs / p_bd_months / p_bd_day), each with its own constraint.
#   - survey sheet: turn the "p_sample_collected" question from
#     select_one into select_multiple.
#   - choices sheet: give the "samples" choice list dot-separated names
#     (so they are valid tokens inside a select_multiple answer).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# ---------------------------------------------------------------------------
# 1. Make room for the new birth-data group: insert 4 blank rows right after
#    row 7 (the old "date / p_birth_date" row), pushing the old rows 8-23
#    down to rows 12-27.
# ---------------------------------------------------------------------------
$survey.Range("A8:A11").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Row 7 becomes the group opener (was the single date question).
# ---------------------------------------------------------------------------
$survey.Range("A7").Value = "begin group"
$survey.Range("B7").Value = "birth_data"
$survey.Range("C7").Value = "1.3. Fecha nacimiento "
$survey.Range("D7").Value = $null
$survey.Range("E7").Value = "field-list"
$survey.Range("F7").Value = $null
$survey.Range("G7").Value = $null
$survey.Range("H7").Value = $null
$survey.Range("I7").Value = $null

# ---------------------------------------------------------------------------
# 3. Rows 8-10: the three new integer sub-questions.
# ---------------------------------------------------------------------------
$survey.Range("A8").Value = "integer"
$survey.Range("B8").Value = "p_bd_yrs"
$survey.Range("C8").Value = "1.3.1. Año de nacimiento"
$survey.Range("F8").Value = '${p_bd_yrs} > 1900 and ${p_bd_yrs} < 2005'
$survey.Range("G8").Value = "El año es incorrecto"
$survey.Range("I8").Value = "yes"

$survey.Range("A9").Value = "integer"
$survey.Range("B9").Value = "p_bd_months"
$survey.Range("C9").Value = "1.3.2. Mes de nacimiento"
$survey.Range("F9").Value = '${p_bd_months} <= 12 and ${p_bd_months} > 0'
$survey.Range("G9").Value = "El mes es incorrecto"
$survey.Range("I9").Value = "yes"

$survey.Range("A10").Value = "integer"
$survey.Range("B10").Value = "p_bd_day"
$survey.Range("C10").Value = "1.3.3. Día de nacimiento"
$survey.Range("F10").Value = '${p_bd_day} > 0 and ${p_bd_day} <= 31'
$survey.Range("G10").Value = "Los días son incorrectos"
$survey.Range("I10").Value = "yes"

# ---------------------------------------------------------------------------
# 4. Row 11: the group closer.
# ---------------------------------------------------------------------------
$survey.Range("A11").Value = "end group"

# ---------------------------------------------------------------------------
# 5. The "Muestras recogidas" question (old row 22, now row 26 after the
#    shift) turns into a select_multiple.
# ---------------------------------------------------------------------------
$survey.Range("I2").Copy()
$survey.Range("A26").PasteSpecial(-4122)
$survey.Range("A26").Value = "select_multiple samples"

# ---------------------------------------------------------------------------
# 6. choices sheet: "samples" list option names become dot-separated so they
#    are usable as select_multiple tokens (labels stay the same).
# ---------------------------------------------------------------------------
$choices.Range("B36").Value = "Papel.whatman"
$choices.Range("B37").Value = "Gota.gruesa"

# ---------------------------------------------------------------------------
# 7. Column sizing: columns A and F grew wider to fit the new labels /
#    constraint formulas.
# ---------------------------------------------------------------------------
$survey.Columns.Item(1).ColumnWidth = 30.41015625
$survey.Columns.Item(6).ColumnWidth = 42.46875

# ---------------------------------------------------------------------------
# 8. Refresh the remembered selection on the survey sheet (bottom-left,
#    unfrozen pane now parks on A2 instead of the old C22).
# ---------------------------------------------------------------------------
$survey.Range("A2").Select()

# ---------------------------------------------------------------------------
# 9. choices sheet view had scrolled down with D45 selected.
# ---------------------------------------------------------------------------
$choices.Range("D45").Select()
$survey.Activate()
